# "add verifying table demo"
# Insert a "Type" column at the front and a "Priority" column before "Created",
# plus append a "Due Date" column at the end; populate the new cells, fix up
# the hyperlinks/date formatting/drawing anchors that the column inserts
# otherwise leave stale.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Add the "Due Date" column (H) first - matches the shared-string order
#    produced by the original authoring session (Due Date is allocated before
#    Type/Bug/Priority/Medium).
# ---------------------------------------------------------------------------
$ws.Range("F1").Value = "Due Date"

# ---------------------------------------------------------------------------
# 2. Insert the new "Type" column before column A. Existing A:E shifts to B:F.
# ---------------------------------------------------------------------------
$ws.Columns("A:A").Insert()
$ws.Range("A1").Value = "Type"
$ws.Range("A2").Value = "Bug"
$ws.Range("A3").Value = "Bug"
$ws.Range("A4").Value = "Bug"
$ws.Range("A5").Value = "Bug"

# ---------------------------------------------------------------------------
# 3. Insert the new "Priority" column before the (now shifted) "Created"
#    column, i.e. before column E.
# ---------------------------------------------------------------------------
$ws.Columns("E:E").Insert()
$ws.Range("E1").Value = "Priority"
$ws.Range("E2").Value = "Medium"
$ws.Range("E3").Value = "Medium"
$ws.Range("E4").Value = "Medium"
$ws.Range("E5").Value = "Medium"

# ---------------------------------------------------------------------------
# 4. Re-point the hyperlinks (column insert does not move them automatically).
#    Leave TextToDisplay unspecified so the existing cell text (shifted into
#    place by the column inserts above) is preserved untouched.
# ---------------------------------------------------------------------------
$ws.Range("A1").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("B2"), "https://katalon.atlassian.net/browse/KD-24272") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://katalon.atlassian.net/browse/KD-24272") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "https://katalon.atlassian.net/secure/ViewProfile.jspa?name=demo") | Out-Null

$ws.Hyperlinks.Add($ws.Range("B3"), "https://katalon.atlassian.net/browse/KD-24229") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://katalon.atlassian.net/browse/KD-24229") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "https://katalon.atlassian.net/secure/ViewProfile.jspa?name=demo") | Out-Null

$ws.Hyperlinks.Add($ws.Range("B4"), "https://katalon.atlassian.net/browse/KD-24197") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "https://katalon.atlassian.net/browse/KD-24197") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D4"), "https://katalon.atlassian.net/secure/ViewProfile.jspa?name=demo") | Out-Null

$ws.Hyperlinks.Add($ws.Range("B5"), "https://katalon.atlassian.net/browse/KD-24191") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C5"), "https://katalon.atlassian.net/browse/KD-24191") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D5"), "https://katalon.atlassian.net/secure/ViewProfile.jspa?name=demo") | Out-Null

# ---------------------------------------------------------------------------
# 5. New Due Date column cells (H2:H5) - blank but date-formatted, same style
#    as Created/Updated.
# ---------------------------------------------------------------------------
$ws.Range("H2:H5").NumberFormat = "dd/mmm/yy"
$ws.Range("F2:G5").NumberFormat = "dd/mmm/yy"

# ---------------------------------------------------------------------------
# 6. Re-anchor the floating "Bug"/"Medium" shapes: they used to sit over
#    column A (index 0) and column D (index 3); the column inserts pushed the
#    real columns to B and F, so nudge each shape's left edge onto the new
#    column boundary (keeps the same col-offset of 0 the originals had).
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $ws.Shapes.Count; $i++) {
    $sh = $ws.Shapes.Item($i)
    if ($i % 2 -eq 1) {
        $sh.Left = $ws.Columns.Item(2).Left
    } else {
        $sh.Left = $ws.Columns.Item(6).Left
    }
}

# ---------------------------------------------------------------------------
# 7. Selection / view bookkeeping to match the saved state.
# ---------------------------------------------------------------------------
$ws.Range("H7").Select()
